# Applies the S2_Table supplement edit:
#  1. Move the "_GoBack" bookmark from its position (between the "2" and
#     " Table." runs of the "S2 Table. Satellite Imagery." paragraph) to the
#     very start of that same paragraph (before the "S" run).
#  2. Delete the leading "Supplement S2 Table" heading paragraph and the
#     blank paragraph that follows it, leaving "S2 Table. Satellite
#     Imagery." as the first paragraph in the document.

$d = $word.ActiveDocument

# Locate the paragraph that currently contains "S2 Table. Satellite Imagery."
# (it's the 3rd paragraph in the original layout) and re-anchor the
# "_GoBack" bookmark to its start.
$targetPara = $d.Paragraphs.Item(3)
$targetStart = $targetPara.Range.Start

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
$anchor = $d.Range($targetStart, $targetStart)
$d.Bookmarks.Add("_GoBack", $anchor)

# Remove the first two paragraphs ("Supplement S2 Table" heading and the
# blank paragraph beneath it).
$firstPara = $d.Paragraphs.Item(1)
$secondPara = $d.Paragraphs.Item(2)
$removeRange = $d.Range($firstPara.Range.Start, $secondPara.Range.End)
$removeRange.Delete()
